$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.6

# Row 5
$ws.Range("S5").Value = 1.62

# Row 8
$ws.Range("J8").Value = 4.25
$ws.Range("L8").Value = 2.4
$ws.Range("O8").Value = 1.25
$ws.Range("P8").Value = 3.25
$ws.Range("V8").Value = 1.98
$ws.Range("W8").Value = 11.75
$ws.Range("Y8").Value = 13
$ws.Range("AD8").Value = 6.8
$ws.Range("AG8").Value = 8
$ws.Range("AJ8").Value = 16
$ws.Range("AO8").Value = 22
$ws.Range("AP8").Value = 27
$ws.Range("AQ8").Value = 120
$ws.Range("AS8").Value = 350
$ws.Range("AT8").Value = 2.62
$ws.Range("AU8").Value = 6.9
$ws.Range("AV8").Value = 60
$ws.Range("AW8").Value = 3.75
$ws.Range("AX8").Value = 9
$ws.Range("BA8").Value = 55

# Row 23
$ws.Range("G23").Value = 3
$ws.Range("I23").Value = 2.15
$ws.Range("J23").Value = 3.4
$ws.Range("L23").Value = 2.75
$ws.Range("W23").Value = 15
$ws.Range("X23").Value = 19
$ws.Range("Z23").Value = 34
$ws.Range("AA23").Value = 21
$ws.Range("AH23").Value = 13
$ws.Range("AI23").Value = 9
$ws.Range("AJ23").Value = 21
$ws.Range("AK23").Value = 15
$ws.Range("AN23").Value = 5.5
$ws.Range("AW23").Value = 4.5
$ws.Range("AX23").Value = 11
$ws.Range("AZ23").Value = 34
$ws.Range("BA23").Value = 41
$ws.Range("BD23").Value = 176

# Row 25
$ws.Range("I25").Value = 2.8
$ws.Range("P25").Value = 4.35
$ws.Range("R25").Value = 2.32
$ws.Range("T25").Value = 3.3
$ws.Range("V25").Value = 2.42
$ws.Range("X25").Value = 12.5
$ws.Range("Z25").Value = 21
$ws.Range("AD25").Value = 7.8
$ws.Range("AG25").Value = 13
$ws.Range("AH25").Value = 17.5
$ws.Range("AL25").Value = 23
$ws.Range("AM25").Value = 200
$ws.Range("AT25").Value = 3.3

# Row 26
$ws.Range("G26").Value = 2.52
$ws.Range("H26").Value = 2.72
$ws.Range("I26").Value = 3.1
$ws.Range("J26").Value = 3.15
$ws.Range("L26").Value = 3.8
$ws.Range("M26").Value = 1.1
$ws.Range("N26").Value = 6.78
$ws.Range("O26").Value = 1.47
$ws.Range("P26").Value = 2.32
$ws.Range("Q26").Value = 2.37
$ws.Range("R26").Value = 1.45
$ws.Range("X26").Value = 11.5
$ws.Range("Y26").Value = 9.75
$ws.Range("Z26").Value = 29
$ws.Range("AA26").Value = 24
$ws.Range("AB26").Value = 40
$ws.Range("AC26").Value = 6.2
$ws.Range("AD26").Value = 5.4
$ws.Range("AG26").Value = 7.1
$ws.Range("AH26").Value = 14.5
$ws.Range("AI26").Value = 11.5
$ws.Range("AJ26").Value = 45
$ws.Range("AK26").Value = 35
$ws.Range("AL26").Value = 50
$ws.Range("AN26").Value = 4.2
$ws.Range("AO26").Value = 14
$ws.Range("AP26").Value = 23
$ws.Range("AQ26").Value = 65
$ws.Range("AR26").Value = 110
$ws.Range("AU26").Value = 7.1
$ws.Range("AV26").Value = 75
$ws.Range("AW26").Value = 4.8
$ws.Range("AX26").Value = 18.5
$ws.Range("AY26").Value = 28
$ws.Range("AZ26").Value = 100
$ws.Range("BA26").Value = 150
$ws.Range("BB26").Value = 450

# Row 27
$ws.Range("G27").Value = 2.3
$ws.Range("I27").Value = 3.6
$ws.Range("J27").Value = 2.9
$ws.Range("K27").Value = 1.87
$ws.Range("L27").Value = 4.2
$ws.Range("N27").Value = 6.15
$ws.Range("Q27").Value = 2.35
$ws.Range("S27").Value = 1.5
$ws.Range("T27").Value = 2.25
$ws.Range("W27").Value = 6.1
$ws.Range("X27").Value = 10.25
$ws.Range("Y27").Value = 9
$ws.Range("AC27").Value = 6.2
$ws.Range("AG27").Value = 8
$ws.Range("AH27").Value = 18.5
$ws.Range("AJ27").Value = 60
$ws.Range("AN27").Value = 4
$ws.Range("AP27").Value = 21
$ws.Range("AT27").Value = 2.22
$ws.Range("AU27").Value = 6.8
$ws.Range("AV27").Value = 65
$ws.Range("AW27").Value = 5.3
$ws.Range("AX27").Value = 22
